$d = $word.ActiveDocument

function Set-CellText($rowIndex, $newText) {
    $tbl = $d.Tables.Item(1)
    $cell = $tbl.Cell($rowIndex, 1)
    $cell.Range.Text = $newText
}

Set-CellText 1  "0M"
Set-CellText 2  "0M"
Set-CellText 3  "0M"
Set-CellText 4  "3513"
Set-CellText 6  "0.47732"
Set-CellText 7  "0.06674"
Set-CellText 8  "0.00801"
Set-CellText 9  "0.41524"
Set-CellText 10 "0.43684"
Set-CellText 11 "0.44043"
Set-CellText 12 "5.94489"

Set-CellText 44 "99.53"
Set-CellText 45 "5.94"
Set-CellText 46 "1261"
